# Updated my sections for the Itus Review
#
# Reconstructs the comment list in the document so that:
#   - "Sections 1-3 (Lloyd)"  -> "Sections 1-3 (Beaver)"
#   - a new "2.2 ..." review comment replaces the old stray "5.2.4 ..." bullet
#     that used to sit right after "Sections 4-6", and carries the _GoBack
#     bookmark that used to live on the last paragraph
#   - a new "3.1.4 ..." review comment is added
#   - "Sections  4-6 (Jeff)" -> "Sections  4-6 (Bandit)" and moves down,
#     now followed by the (re-added) "5.2.4 ..." bullet and then "6.2.1 ..."
#   - "Sections 7-9 (Tyler)" -> "Sections 7-9 (Penny)"
#   - "Sections 10-12 (CM)"  -> "Sections 10-12 (Dawg)"
#
# The whole body is rebuilt in one shot via Range.InsertXML so that the
# exact run layout (including the split runs around the replaced names and
# the proofErr/bookmark placement) matches what Word produced for this edit.

$d = $word.ActiveDocument

$body = $d.Content

$xml = @'
<w:p w:rsidR="00F644AB" w:rsidRDefault="00C02A68" w:rsidP="00F94E67"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>COMMENTS:</w:t></w:r></w:p><w:p w:rsidR="00F94E67" w:rsidRDefault="00F94E67" w:rsidP="00F94E67"><w:r><w:t>Overall:</w:t></w:r></w:p><w:p w:rsidR="00C02A68" w:rsidRDefault="00C02A68" w:rsidP="00C02A68"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve"> There are too many document revision history entries. Only working versions of the document should be entered.</w:t></w:r></w:p><w:p w:rsidR="00F94E67" w:rsidRDefault="00F94E67" w:rsidP="00C02A68"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>The formatting makes it a little hard to differentiate between individual requirements. If the sub parts of each requirement (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> x.x.1-x.x.5) were indented or smaller text it would improve flow and readability of the document.</w:t></w:r></w:p><w:p w:rsidR="00C02A68" w:rsidRDefault="00F94E67" w:rsidP="00ED366C"><w:r><w:t>Sections 1-3 (</w:t></w:r><w:r><w:t>Beaver</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p><w:p w:rsidR="00ED366C" w:rsidRDefault="00ED366C" w:rsidP="00ED366C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>2.2</w:t></w:r><w:r><w:t xml:space="preserve"> What about changing </w:t></w:r><w:r><w:t>frequencies</w:t></w:r><w:r><w:t xml:space="preserve">? I’ve been paintballing and many sites have parties with 10-15 </w:t></w:r><w:r><w:t>paintball</w:t></w:r><w:r><w:t xml:space="preserve"> courses. And they just rotate you through the different courses for the whole day. I think </w:t></w:r><w:r><w:t>changing frequencies would be something you want the user to be able to set.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p w:rsidR="00ED366C" w:rsidRDefault="00ED366C" w:rsidP="00ED366C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">3.1.4 – What about FCC requirements; </w:t></w:r><w:r><w:t>is</w:t></w:r><w:r><w:t xml:space="preserve"> that an issue? For example you can transmit radio on 92.5 but if you transmit more than 100 ft. you now have a pirate radio signal.</w:t></w:r></w:p><w:p w:rsidR="00ED366C" w:rsidRDefault="00ED366C" w:rsidP="00ED366C"/><w:p w:rsidR="00ED366C" w:rsidRDefault="00ED366C" w:rsidP="00ED366C"><w:proofErr w:type="gramStart"/><w:r><w:t>Sections  4</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>-6 (</w:t></w:r><w:r><w:t>Bandit</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p><w:p w:rsidR="00ED366C" w:rsidRDefault="00ED366C" w:rsidP="00ED366C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>5.2.4 – Human reaction time not really a standard, can very broadly.</w:t></w:r></w:p><w:p w:rsidR="00ED366C" w:rsidRDefault="00ED366C" w:rsidP="00F94E67"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>6.2.1 – No obvious reason to specify ‘green’ light</w:t></w:r></w:p><w:p w:rsidR="00F94E67" w:rsidRDefault="00F94E67" w:rsidP="00F94E67"/><w:p w:rsidR="00F94E67" w:rsidRDefault="00F94E67" w:rsidP="00F94E67"><w:r><w:t>Sections 7-9 (</w:t></w:r><w:r><w:t>Penny</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p><w:p w:rsidR="00F94E67" w:rsidRDefault="00F94E67" w:rsidP="00F94E67"/><w:p w:rsidR="00F94E67" w:rsidRDefault="00F94E67" w:rsidP="00F94E67"><w:r><w:t>Sections 10-12 (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Dawg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>
'@

$body.InsertXML($xml) | Out-Null
